# repull data, push all data, mean calculation
# Update the dSF column (F) with refreshed/repulled values for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = 0
    6  = -6
    7  = 1
    10 = -4
    11 = -4
    12 = -3
    13 = 0
    17 = 1
    19 = 3
    23 = 0
    30 = -2
    34 = -1
    36 = -1
    38 = 5
    40 = 5
    48 = -4
    56 = -4
    57 = -11
    58 = 2
    60 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
